$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before column L; this shifts old L..T to M..U
$ws.Columns("L:L").Insert()

# New column headers for the "Exceptional items" column
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"

# New column values (computed as P/l before tax - P/l before exceptional items & tax)
$ws.Range("L3").Value = $null
$ws.Range("L4").Value = $null
$ws.Range("L5").Value = $null
$ws.Range("L6").Value = $null
$ws.Range("L7").Value = -11.26
$ws.Range("L8").Value = $null
$ws.Range("L9").Value = $null
$ws.Range("L10").Value = -6.09
$ws.Range("L11").Value = $null
$ws.Range("L12").Value = 6.09
$ws.Range("L13").Value = $null
$ws.Range("L14").Value = $null
$ws.Range("L15").Value = $null
$ws.Range("L16").Value = $null
$ws.Range("L17").Value = $null
$ws.Range("L18").Value = 1.07
$ws.Range("L19").Value = $null
$ws.Range("L20").Value = $null
$ws.Range("L21").Value = $null
$ws.Range("L22").Value = 0.36
